$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I1").Value = "Then_Question"
$ws.Range("J1").Value = "Else_Question"

$ws.Range("J1").Select()
